# Update coworker notebook: refresh monthly data points for
# 中国10年期国债收益率Trend/F0.2 (column B "实际值" actuals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 1.59
$ws.Range("B6").Value = 1.68
$ws.Range("B11").Value = 1.89
$ws.Range("B19").Value = 2.21
